$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.620.51"
$ws.Range("E2").Value = "  -1.69%  "
$ws.Range("D3").Value = "3.629.20"
$ws.Range("E3").Value = "  -0.22%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'583.79"
$ws.Range("E5").Value = "  -1.49%  "
$ws.Range("D6").Value = "'175.33"
$ws.Range("E6").Value = "  -2.82%  "
$ws.Range("D7").Value = "'0.632"
$ws.Range("E7").Value = "  +4.28%  "
$ws.Range("D8").Value = "3.620.13"
$ws.Range("E8").Value = "  -0.23%  "
$ws.Range("D9").Value = "'1.00"
$ws.Range("E9").Value = "  -0.14%  "
$ws.Range("D10").Value = "'0.194"
$ws.Range("E10").Value = "  -4.45%  "
$ws.Range("D11").Value = "'6.64"
$ws.Range("E11").Value = "  +14.66%  "
$ws.Range("D12").Value = "'0.617"
$ws.Range("E12").Value = "  +1.95%  "
$ws.Range("D13").Value = "'48.42"
$ws.Range("E13").Value = "  -3.01%  "
$ws.Range("D14").Value = "'0.0000281"
$ws.Range("E14").Value = "  -1.69%  "
$ws.Range("D15").Value = "'682.64"
$ws.Range("E15").Value = "  -1.78%  "
$ws.Range("D16").Value = "4.218.92"
$ws.Range("E16").Value = "  -0.53%  "
$ws.Range("D17").Value = "'9.02"
$ws.Range("E17").Value = "  +0.36%  "
$ws.Range("D18").Value = "3.623.48"
$ws.Range("E18").Value = "  -0.79%  "
$ws.Range("D19").Value = "70.703.12"
$ws.Range("E19").Value = "  -1.88%  "
$ws.Range("E20").Value = "  -0.58%  "
$ws.Range("D21").Value = "'17.78"
$ws.Range("E21").Value = "  -3.49%  "
$ws.Range("E22").Value = "  -0.95%  "
$ws.Range("D23").Value = "'0.937"
$ws.Range("E23").Value = "  +0.32%  "
$ws.Range("D24").Value = "'17.04"
$ws.Range("E24").Value = "  -4.67%  "
$ws.Range("D25").Value = "'99.90"
$ws.Range("E25").Value = "  -3.64%  "
$ws.Range("E26").Value = "  -3.11%  "
$ws.Range("D27").Value = "'2.77"
$ws.Range("E27").Value = "  -3.12%  "
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("D29").Value = "'9.78"
$ws.Range("E29").Value = "  -1.86%  "
$ws.Range("D30").Value = "'34.42"
$ws.Range("E30").Value = "  -1.79%  "
$ws.Range("D31").Value = "'9.14"
$ws.Range("E31").Value = "  +0.23%  "
$ws.Range("D32").Value = "'3.29"
$ws.Range("E32").Value = "  -4.11%  "
$ws.Range("D33").Value = "'7.47"
$ws.Range("E33").Value = "  +2.21%  "
$ws.Range("E34").Value = "  -4.23%  "
$ws.Range("D35").Value = "'3.94"
$ws.Range("E35").Value = "  -5.36%  "
$ws.Range("D36").Value = "'570.98"
$ws.Range("E36").Value = "  -1.88%  "
$ws.Range("E37").Value = "  -2.22%  "
$ws.Range("E38").Value = "  -1.73%  "
$ws.Range("E39").Value = "  -1.52%  "
$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("E41").Value = "  -0.70%  "
$ws.Range("D42").Value = "'0.0450"
$ws.Range("E42").Value = "  -2.51%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "3.534.43"
$ws.Range("E43").Value = "  -3.43%  "
$ws.Range("B44").Value = "Kaspa"
$ws.Range("C44").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D44").Value = "'0.140"
$ws.Range("E44").Value = "  -2.42%  "
$ws.Range("D45").Value = "'34.24"
$ws.Range("E45").Value = "  -4.29%  "
$ws.Range("D46").Value = "0.0₃0727"
$ws.Range("E46").Value = "  -5.17%  "
$ws.Range("D47").Value = "'3.01"
$ws.Range("E47").Value = "  +6.03%  "
$ws.Range("E48").Value = "  -3.71%  "
$ws.Range("E49").Value = "  +1.84%  "
$ws.Range("D50").Value = "'137.30"
$ws.Range("E50").Value = "  +3.70%  "
$ws.Range("D51").Value = "'2.87"
$ws.Range("E51").Value = "  -3.94%  "
